$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.070536666666667
$ws.Range("H2").Value = 9.21161
$ws.Range("I2").Value = 0.2245190988242715
$ws.Range("J2").Value = 0.2245190988242715
$ws.Range("O2").Value = 0.6892208999344782
$ws.Range("P2").Value = 0.6892208999344781
$ws.Range("Q2").Value = 0.83330987543
$ws.Range("R2").Value = 7.49978887887
$ws.Range("S2").Value = 0.1547432553441425
$ws.Range("T2").Value = 0.1547432553441424

$ws.Range("G3").Value = 3.070536666666667
$ws.Range("H3").Value = 9.21161
$ws.Range("I3").Value = 0.2245190988242715
$ws.Range("J3").Value = 0.2245190988242715
$ws.Range("M3").Value = 0.122373
$ws.Range("N3").Value = 0.367119
$ws.Range("O3").Value = 0.3107791000655218
$ws.Range("P3").Value = 0.3107791000655218
$ws.Range("Q3").Value = 0.37575078351
$ws.Range("R3").Value = 3.38175705159
$ws.Range("S3").Value = 0.06977584348012905
$ws.Range("T3").Value = 0.06977584348012904

$ws.Range("I4").Value = 0.2851553493233187
$ws.Range("J4").Value = 0.2851553493233187
$ws.Range("O4").Value = 0.6892208999344782
$ws.Range("P4").Value = 0.6892208999344781
$ws.Range("S4").Value = 0.1965350264817482
$ws.Range("T4").Value = 0.1965350264817482

$ws.Range("I5").Value = 0.2851553493233187
$ws.Range("J5").Value = 0.2851553493233187
$ws.Range("M5").Value = 0.122373
$ws.Range("N5").Value = 0.367119
$ws.Range("O5").Value = 0.3107791000655218
$ws.Range("P5").Value = 0.3107791000655218
$ws.Range("Q5").Value = 0.477230429355
$ws.Range("R5").Value = 4.295073864194999
$ws.Range("S5").Value = 0.08862032284157048
$ws.Range("T5").Value = 0.08862032284157048

$ws.Range("G6").Value = 3.959514
$ws.Range("H6").Value = 11.878542
$ws.Range("I6").Value = 0.2895215434854775
$ws.Range("J6").Value = 0.2895215434854775
$ws.Range("O6").Value = 0.6892208999344782
$ws.Range("P6").Value = 0.6892208999344781
$ws.Range("Q6").Value = 1.074568544946
$ws.Range("R6").Value = 9.671116904514
$ws.Range("S6").Value = 0.19954429875148
$ws.Range("T6").Value = 0.1995442987514799

$ws.Range("G7").Value = 3.959514
$ws.Range("H7").Value = 11.878542
$ws.Range("I7").Value = 0.2895215434854775
$ws.Range("J7").Value = 0.2895215434854775
$ws.Range("M7").Value = 0.122373
$ws.Range("N7").Value = 0.367119
$ws.Range("O7").Value = 0.3107791000655218
$ws.Range("P7").Value = 0.3107791000655218
$ws.Range("Q7").Value = 0.484537606722
$ws.Range("R7").Value = 4.360838460497999
$ws.Range("S7").Value = 0.08997724473399753
$ws.Range("T7").Value = 0.08997724473399753

$ws.Range("G8").Value = 2.746207666666667
$ws.Range("H8").Value = 8.238623
$ws.Range("I8").Value = 0.2008040083669322
$ws.Range("J8").Value = 0.2008040083669322
$ws.Range("O8").Value = 0.6892208999344782
$ws.Range("P8").Value = 0.6892208999344781
$ws.Range("Q8").Value = 0.745290552449
$ws.Range("R8").Value = 6.707614972041
$ws.Range("S8").Value = 0.1383983193571075
$ws.Range("T8").Value = 0.1383983193571075

$ws.Range("G9").Value = 2.746207666666667
$ws.Range("H9").Value = 8.238623
$ws.Range("I9").Value = 0.2008040083669322
$ws.Range("J9").Value = 0.2008040083669322
$ws.Range("M9").Value = 0.122373
$ws.Range("N9").Value = 0.367119
$ws.Range("O9").Value = 0.3107791000655218
$ws.Range("P9").Value = 0.3107791000655218
$ws.Range("Q9").Value = 0.336061670793
$ws.Range("R9").Value = 3.024555037137
$ws.Range("S9").Value = 0.06240568900982469
$ws.Range("T9").Value = 0.06240568900982469

